$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 112, shifting existing rows 112:247 down to 113:248
$ws.Rows("112").Insert()

# Populate the newly inserted row 112 with its values
$ws.Range("A112").Value = 8
$ws.Range("B112").Value = "Terminal La Palmera de La Serena"
$ws.Range("C112").Value = "Coquimbo"
$ws.Range("D112").Value = 44546
$ws.Range("E112").Value = 4
$ws.Range("F112").Value = 100114013
$ws.Range("G112").Value = "Zanahoria"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 600
$ws.Range("K112").Value = 6000
$ws.Range("L112").Value = 7000
$ws.Range("M112").Value = 6500
$ws.Range("N112").Value = "$/saco 20 kilos"
$ws.Range("O112").Value = "Provincia del Elquí"
$ws.Range("P112").Value = 325
$ws.Range("Q112").Value = 20
$ws.Range("R112").Value = "Hortaliza"
